# Update the "Förändrad" (changed) date column C for rows 2-6 on the
# single worksheet "Avverkningsanmälningar" (Översikt SÖDERHAMN).
# Old value: 45233 (2023-11-03) -> New value: 45243 (2023-11-13)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..6) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
